$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.2711976666666667
$ws.Range("H2").Value = 0.813593
$ws.Range("I2").Value = 0.2043600193410237
$ws.Range("J2").Value = 0.2043600193410237
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.709791333333333
$ws.Range("N2").Value = 11.129374
$ws.Range("O2").Value = 0.4283284425582907
$ws.Range("P2").Value = 0.4283284425582907
$ws.Range("Q2").Value = 1.006086753420222
$ws.Range("R2").Value = 9.054780780782
$ws.Range("S2").Value = 0.08753320880552283
$ws.Range("T2").Value = 0.08753320880552286

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.2711976666666667
$ws.Range("H3").Value = 0.813593
$ws.Range("I3").Value = 0.2043600193410237
$ws.Range("J3").Value = 0.2043600193410237
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.599001333333334
$ws.Range("N3").Value = 10.797004
$ws.Range("O3").Value = 0.4155367505499981
$ws.Range("P3").Value = 0.4155367505499982
$ws.Range("Q3").Value = 0.9760407639302223
$ws.Range("R3").Value = 8.784366875372001
$ws.Range("S3").Value = 0.08491909837930375
$ws.Range("T3").Value = 0.08491909837930377

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.2711976666666667
$ws.Range("H4").Value = 0.813593
$ws.Range("I4").Value = 0.2043600193410237
$ws.Range("J4").Value = 0.2043600193410237
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.352297666666667
$ws.Range("N4").Value = 4.056893
$ws.Range("O4").Value = 0.1561348068917112
$ws.Range("P4").Value = 0.1561348068917112
$ws.Range("Q4").Value = 0.3667399718387778
$ws.Range("R4").Value = 3.300659746549
$ws.Range("S4").Value = 0.03190771215619709
$ws.Range("T4").Value = 0.0319077121561971

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.5110070000000001
$ws.Range("H5").Value = 1.533021
$ws.Range("I5").Value = 0.3850674737985645
$ws.Range("J5").Value = 0.3850674737985645
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.709791333333333
$ws.Range("N5").Value = 11.129374
$ws.Range("O5").Value = 0.4283284425582907
$ws.Range("P5").Value = 0.4283284425582907
$ws.Range("Q5").Value = 1.895729339872667
$ws.Range("R5").Value = 17.061564058854
$ws.Range("S5").Value = 0.1649353513319946
$ws.Range("T5").Value = 0.1649353513319946

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.5110070000000001
$ws.Range("H6").Value = 1.533021
$ws.Range("I6").Value = 0.3850674737985645
$ws.Range("J6").Value = 0.3850674737985645
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.599001333333334
$ws.Range("N6").Value = 10.797004
$ws.Range("O6").Value = 0.4155367505499981
$ws.Range("P6").Value = 0.4155367505499982
$ws.Range("Q6").Value = 1.839114874342667
$ws.Range("R6").Value = 16.552033869084
$ws.Range("S6").Value = 0.1600096868047521
$ws.Range("T6").Value = 0.1600096868047521

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.5110070000000001
$ws.Range("H7").Value = 1.533021
$ws.Range("I7").Value = 0.3850674737985645
$ws.Range("J7").Value = 0.3850674737985645
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.352297666666667
$ws.Range("N7").Value = 4.056893
$ws.Range("O7").Value = 0.1561348068917112
$ws.Range("P7").Value = 0.1561348068917112
$ws.Range("Q7").Value = 0.6910335737503335
$ws.Range("R7").Value = 6.219302163753
$ws.Range("S7").Value = 0.06012243566181794
$ws.Range("T7").Value = 0.06012243566181794

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.181472
$ws.Range("H8").Value = 0.544416
$ws.Range("I8").Value = 0.1367475682430438
$ws.Range("J8").Value = 0.1367475682430438
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.709791333333333
$ws.Range("N8").Value = 11.129374
$ws.Range("O8").Value = 0.4283284425582907
$ws.Range("P8").Value = 0.4283284425582907
$ws.Range("Q8").Value = 0.6732232528426667
$ws.Range("R8").Value = 6.059009275584001
$ws.Range("S8").Value = 0.05857287292917653
$ws.Range("T8").Value = 0.05857287292917654

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.181472
$ws.Range("H9").Value = 0.544416
$ws.Range("I9").Value = 0.1367475682430438
$ws.Range("J9").Value = 0.1367475682430438
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.599001333333334
$ws.Range("N9").Value = 10.797004
$ws.Range("O9").Value = 0.4155367505499981
$ws.Range("P9").Value = 0.4155367505499982
$ws.Range("Q9").Value = 0.6531179699626667
$ws.Range("R9").Value = 5.878061729664001
$ws.Range("S9").Value = 0.05682364015332855
$ws.Range("T9").Value = 0.05682364015332855

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.181472
$ws.Range("H10").Value = 0.544416
$ws.Range("I10").Value = 0.1367475682430438
$ws.Range("J10").Value = 0.1367475682430438
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.352297666666667
$ws.Range("N10").Value = 4.056893
$ws.Range("O10").Value = 0.1561348068917112
$ws.Range("P10").Value = 0.1561348068917112
$ws.Range("Q10").Value = 0.2454041621653333
$ws.Range("R10").Value = 2.208637459488
$ws.Range("S10").Value = 0.02135105516053875
$ws.Range("T10").Value = 0.02135105516053875

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.08639666666666668
$ws.Range("H11").Value = 0.25919
$ws.Range("I11").Value = 0.06510389520681709
$ws.Range("J11").Value = 0.06510389520681709
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.709791333333333
$ws.Range("N11").Value = 11.129374
$ws.Range("O11").Value = 0.4283284425582907
$ws.Range("P11").Value = 0.4283284425582907
$ws.Range("Q11").Value = 0.3205136052288889
$ws.Range("R11").Value = 2.88462244706
$ws.Range("S11").Value = 0.02788585003841413
$ws.Range("T11").Value = 0.02788585003841413

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.08639666666666668
$ws.Range("H12").Value = 0.25919
$ws.Range("I12").Value = 0.06510389520681709
$ws.Range("J12").Value = 0.06510389520681709
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.599001333333334
$ws.Range("N12").Value = 10.797004
$ws.Range("O12").Value = 0.4155367505499981
$ws.Range("P12").Value = 0.4155367505499982
$ws.Range("Q12").Value = 0.310941718528889
$ws.Range("R12").Value = 2.798475466760001
$ws.Range("S12").Value = 0.02705306106238837
$ws.Range("T12").Value = 0.02705306106238837

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.08639666666666668
$ws.Range("H13").Value = 0.25919
$ws.Range("I13").Value = 0.06510389520681709
$ws.Range("J13").Value = 0.06510389520681709
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.352297666666667
$ws.Range("N13").Value = 4.056893
$ws.Range("O13").Value = 0.1561348068917112
$ws.Range("P13").Value = 0.1561348068917112
$ws.Range("Q13").Value = 0.1168340107411111
$ws.Range("R13").Value = 1.05150609667
$ws.Range("S13").Value = 0.01016498410601459
$ws.Range("T13").Value = 0.01016498410601459

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.276985
$ws.Range("H14").Value = 0.830955
$ws.Range("I14").Value = 0.2087210434105509
$ws.Range("J14").Value = 0.2087210434105509
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.709791333333333
$ws.Range("N14").Value = 11.129374
$ws.Range("O14").Value = 0.4283284425582907
$ws.Range("P14").Value = 0.4283284425582907
$ws.Range("Q14").Value = 1.027556552463333
$ws.Range("R14").Value = 9.24800897217
$ws.Range("S14").Value = 0.08940115945318262
$ws.Range("T14").Value = 0.08940115945318267

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.276985
$ws.Range("H15").Value = 0.830955
$ws.Range("I15").Value = 0.2087210434105509
$ws.Range("J15").Value = 0.2087210434105509
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.599001333333334
$ws.Range("N15").Value = 10.797004
$ws.Range("O15").Value = 0.4155367505499981
$ws.Range("P15").Value = 0.4155367505499982
$ws.Range("Q15").Value = 0.9968693843133334
$ws.Range("R15").Value = 8.97182445882
$ws.Range("S15").Value = 0.0867312641502254
$ws.Range("T15").Value = 0.08673126415022543

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.276985
$ws.Range("H16").Value = 0.830955
$ws.Range("I16").Value = 0.2087210434105509
$ws.Range("J16").Value = 0.2087210434105509
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.352297666666667
$ws.Range("N16").Value = 4.056893
$ws.Range("O16").Value = 0.1561348068917112
$ws.Range("P16").Value = 0.1561348068917112
$ws.Range("Q16").Value = 0.3667399718387778
$ws.Range("R16").Value = 3.371095522815
$ws.Range("S16").Value = 0.03258861980714283
$ws.Range("T16").Value = 0.03258861980714284

